$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.993.33"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "1.990.56"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'242.19"
$ws.Range("E5").Value = "  -6.23%  "
$ws.Range("D6").Value = "'0.603"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'54.75"
$ws.Range("E8").Value = "  -4.88%  "
$ws.Range("E9").Value = "  -4.11%  "
$ws.Range("D10").Value = "'57.91"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").Value = "'0.0755"
$ws.Range("E11").Value = "  -5.52%  "
$ws.Range("D12").Value = "'0.0978"
$ws.Range("E12").Value = "  -4.53%  "
$ws.Range("D13").Value = "2.280.59"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("E14").Value = "  -5.50%  "
$ws.Range("D15").Value = "'20.60"
$ws.Range("E15").Value = "  -4.04%  "
$ws.Range("E16").Value = "  -8.29%  "
$ws.Range("E17").Value = "  -6.61%  "
$ws.Range("D18").Value = "1.995.55"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").Value = "36.926.82"
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("D20").Value = "'68.09"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("E21").Value = "  -5.76%  "
$ws.Range("D22").Value = "'227.78"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").Value = "'4.98"
$ws.Range("E23").Value = "  -4.81%  "
$ws.Range("E25").Value = "  -9.77%  "
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").Value = "'161.39"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("D28").Value = "'8.64"
$ws.Range("E28").Value = "  -5.71%  "
$ws.Range("D29").Value = "'19.12"
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("D30").Value = "'0.123"
$ws.Range("E30").Value = "  -10.60%  "
$ws.Range("E31").Value = "  -3.55%  "
$ws.Range("D32").Value = "'0.117"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").Value = "'4.42"
$ws.Range("E33").Value = "  -6.96%  "
$ws.Range("D34").Value = "'0.0608"
$ws.Range("E34").Value = "  -9.21%  "
$ws.Range("D35").Value = "'4.21"
$ws.Range("E35").Value = "  -7.55%  "
$ws.Range("E36").Value = "  -6.01%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").Value = "'3.29"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("D40").Value = "'5.20"
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("D41").Value = "'3.11"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("D42").Value = "1.433.17"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  -6.33%  "
$ws.Range("D44").Value = "'1.12"
$ws.Range("E44").Value = "  -6.47%  "
$ws.Range("D45").Value = "'0.0872"
$ws.Range("E45").Value = "  -9.87%  "
$ws.Range("D46").Value = "'88.13"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("D47").Value = "'15.16"
$ws.Range("E47").Value = "  -6.69%  "
$ws.Range("E48").Value = "  -4.95%  "
$ws.Range("D49").Value = "'2.89"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").Value = "'3.66"
$ws.Range("E50").Value = "  +13.81%  "
$ws.Range("D51").Value = "'6.65"
$ws.Range("E51").Value = "  -10.02%  "
